$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a literal text value into a cell without Excel re-typing it as a
# number/date (matches the source data, which stores these as plain text).
function Set-TextValue($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value2 = $text
    $cell.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '30.251.46'

# Row 3
Set-TextValue $ws.Range("D3") '1.878.03'
$ws.Range("E3").Value2 = '  -1.49%  '

# Row 4
$ws.Range("E4").Value2 = '  -0.23%  '

# Row 5
$ws.Range("E5").Value2 = '  -0.59%  '

# Row 6
$ws.Range("E6").Value2 = '  -0.03%  '

# Row 7
Set-TextValue $ws.Range("D7") '0.4852'
$ws.Range("E7").Value2 = '  -0.88%  '

# Row 8
Set-TextValue $ws.Range("D8") '0.2877'
$ws.Range("E8").Value2 = '  -2.71%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.06583'
$ws.Range("E9").Value2 = '  -2.52%  '

# Row 10
Set-TextValue $ws.Range("D10") '1.877.47'
$ws.Range("E10").Value2 = '  -2.34%  '

# Row 11
Set-TextValue $ws.Range("D11") '16.75'
$ws.Range("E11").Value2 = '  -1.43%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.07310'
$ws.Range("E12").Value2 = '  +0.31%  '

# Row 13
Set-TextValue $ws.Range("D13") '5.151'
$ws.Range("E13").Value2 = '  +0.61%  '

# Row 14
Set-TextValue $ws.Range("D14") '87.19'
$ws.Range("E14").Value2 = '  -2.71%  '

# Row 15
Set-TextValue $ws.Range("D15") '0.6550'
$ws.Range("E15").Value2 = '  -1.88%  '

# Row 16
Set-TextValue $ws.Range("D16") '30.224.67'
$ws.Range("E16").Value2 = '  -1.97%  '

# Row 17
$ws.Range("E17").Value2 = '  -0.32%  '

# Row 18
Set-TextValue $ws.Range("D18") '1.0000'
$ws.Range("E18").Value2 = '  +0.04%  '

# Row 19
Set-TextValue $ws.Range("D19") '0.000007731'
$ws.Range("E19").Value2 = '  -2.56%  '

# Row 20
Set-TextValue $ws.Range("D20") '2.131.43'
$ws.Range("E20").Value2 = '  -1.16%  '

# Row 21
Set-TextValue $ws.Range("D21") '5.318'
$ws.Range("E21").Value2 = '  +4.56%  '

# Row 22
Set-TextValue $ws.Range("D22") '1.001'
$ws.Range("E22").Value2 = '  -0.20%  '

# Row 23
Set-TextValue $ws.Range("D23") '194.50'
$ws.Range("E23").Value2 = '  -5.58%  '

# Row 24
Set-TextValue $ws.Range("D24") '6.119'
$ws.Range("E24").Value2 = '  -1.09%  '

# Row 25
Set-TextValue $ws.Range("D25") '9.287'
$ws.Range("E25").Value2 = '  -3.39%  '

# Row 26
Set-TextValue $ws.Range("D26") '161.03'
$ws.Range("E26").Value2 = '  +2.42%  '

# Row 27
Set-TextValue $ws.Range("D27") '18.00'
$ws.Range("E27").Value2 = '  -4.22%  '

# Row 28
Set-TextValue $ws.Range("D28") '1.915'
$ws.Range("E28").Value2 = '  -2.29%  '

# Row 29
Set-TextValue $ws.Range("D29") '1.438'
$ws.Range("E29").Value2 = '  +1.03%  '

# Row 30
Set-TextValue $ws.Range("D30") '4.272'
$ws.Range("E30").Value2 = '  -1.00%  '

# Row 31
Set-TextValue $ws.Range("D31") '0.09121'
$ws.Range("E31").Value2 = '  -0.31%  '

# Row 32
Set-TextValue $ws.Range("D32") '4.033'
$ws.Range("E32").Value2 = '  -0.30%  '

# Row 33
Set-TextValue $ws.Range("D33") '0.05088'
$ws.Range("E33").Value2 = '  -1.54%  '

# Row 34
Set-TextValue $ws.Range("D34") '0.7185'
$ws.Range("E34").Value2 = '  -4.23%  '

# Row 35
Set-TextValue $ws.Range("D35") '1.097'
$ws.Range("E35").Value2 = '  -1.23%  '

# Row 36
$ws.Range("E36").Value2 = '  -0.40%  '

# Row 37
Set-TextValue $ws.Range("D37") '0.01796'
$ws.Range("E37").Value2 = '  -2.20%  '

# Row 38
Set-TextValue $ws.Range("D38") '2.639'
$ws.Range("E38").Value2 = '  -3.26%  '

# Row 39
Set-TextValue $ws.Range("D39") '0.9183'
$ws.Range("E39").Value2 = '  -0.67%  '

# Row 40
Set-TextValue $ws.Range("D40") '2.043'
$ws.Range("E40").Value2 = '  -2.35%  '

# Row 41
Set-TextValue $ws.Range("D41") '106.21'
$ws.Range("E41").Value2 = '  -0.56%  '

# Row 42
Set-TextValue $ws.Range("D42") '0.4277'
$ws.Range("E42").Value2 = '  -4.28%  '

# Row 43
Set-TextValue $ws.Range("D43") '5.801'
$ws.Range("E43").Value2 = '  -1.00%  '

# Row 44
$ws.Range("E44").Value2 = '  -0.58%  '

# Row 45
Set-TextValue $ws.Range("D45") '7.413'
$ws.Range("E45").Value2 = '  -3.74%  '

# Row 46
Set-TextValue $ws.Range("D46") '0.1314'
$ws.Range("E46").Value2 = '  -4.26%  '

# Row 47
Set-TextValue $ws.Range("D47") '65.33'
$ws.Range("E47").Value2 = '  -3.00%  '

# Row 48
Set-TextValue $ws.Range("D48") '8.954'
$ws.Range("E48").Value2 = '  -0.56%  '

# Row 49
$ws.Range("B49").Value2 = 'Elrond'
$ws.Range("C49").Value2 = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
Set-TextValue $ws.Range("D49") '33.87'
$ws.Range("E49").Value2 = '  -3.52%  '

# Row 50
$ws.Range("B50").Value2 = 'Cronos'
$ws.Range("C50").Value2 = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
Set-TextValue $ws.Range("D50") '0.05746'
$ws.Range("E50").Value2 = '  -2.81%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.3819'
$ws.Range("E51").Value2 = '  -6.22%  '
